$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.332.83"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.875.71"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").Value = "1.878.99"
$ws.Range("E13").Value = "  +1.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6950"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.121"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "270.41"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "30.365.09"
$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("E18").Value = "  +1.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007726"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "2.126.47"
$ws.Range("E21").Value = "  +0.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.275"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.219"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.416"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.954"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("E29").Value = "  -1.22%  "

$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.351"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.465"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.073"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04752"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7046"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.724"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.60%  "

$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.806"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.212"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.958"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4184"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8423"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.141"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.198"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "931.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05693"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.88%  "
